$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D2").Value = "2016-02-22 06:35:56"
$wsZh.Range("D3").Value = "2016-02-22 06:35:56"
$wsZh.Range("G2").Value = "2016-02-22 06:37:04"
$wsZh.Range("G3").Value = "2016-02-22 06:37:04"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D2").Value = "2016-02-22 06:36:11"
$wsDe.Range("D3").Value = "2016-02-22 06:36:11"
$wsDe.Range("G2").Value = "2016-02-22 06:37:32"
$wsDe.Range("G3").Value = "2016-02-22 06:37:32"
